# Auto-generated script to update cryptos.xlsx per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $text)
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value2 = $text
    $c.NumberFormat = "General"
    $c.Style = "Normal"
}

Set-CellText $ws 'D2' '30.451.74'
Set-CellText $ws 'E2' '  -0.35%  '

Set-CellText $ws 'D3' '1.933.41'
Set-CellText $ws 'E3' '  +4.36%  '

Set-CellText $ws 'D4' '1.001'
Set-CellText $ws 'E4' '  +0.18%  '

Set-CellText $ws 'D5' '240.40'
Set-CellText $ws 'E5' '  +2.93%  '

Set-CellText $ws 'D6' '1.002'
Set-CellText $ws 'E6' '  +0.21%  '

Set-CellText $ws 'D7' '0.4750'
Set-CellText $ws 'E7' '  +0.42%  '

Set-CellText $ws 'B8' 'OKB'
Set-CellText $ws 'C8' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-CellText $ws 'D8' '44.45'
Set-CellText $ws 'E8' '  +2.51%  '

Set-CellText $ws 'B9' 'Cardano'
Set-CellText $ws 'C9' 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-CellText $ws 'D9' '0.2855'
Set-CellText $ws 'E9' '  +4.21%  '

Set-CellText $ws 'B10' 'Dogecoin'
Set-CellText $ws 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-CellText $ws 'D10' '0.06587'
Set-CellText $ws 'E10' '  +4.26%  '

Set-CellText $ws 'B11' 'Solana'
Set-CellText $ws 'C11' 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-CellText $ws 'D11' '19.31'
Set-CellText $ws 'E11' '  +9.06%  '

Set-CellText $ws 'B12' 'Litecoin'
Set-CellText $ws 'C12' 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-CellText $ws 'D12' '106.82'
Set-CellText $ws 'E12' '  +26.49%  '

Set-CellText $ws 'B13' 'WrappedEther'
Set-CellText $ws 'C13' 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-CellText $ws 'D13' '1.938.70'
Set-CellText $ws 'E13' '  +5.70%  '

Set-CellText $ws 'B14' 'TRON'
Set-CellText $ws 'C14' 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText $ws 'D14' '0.07612'
Set-CellText $ws 'E14' '  +2.24%  '

Set-CellText $ws 'B15' 'Polkadot'
Set-CellText $ws 'C15' 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-CellText $ws 'D15' '5.141'
Set-CellText $ws 'E15' '  +2.89%  '

Set-CellText $ws 'B16' 'Polygon'
Set-CellText $ws 'C16' 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-CellText $ws 'D16' '0.6524'
Set-CellText $ws 'E16' '  +4.29%  '

Set-CellText $ws 'B17' 'BitcoinCash'
Set-CellText $ws 'C17' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText $ws 'D17' '304.41'
Set-CellText $ws 'E17' '  +24.57%  '

Set-CellText $ws 'B18' 'WrappedBTC'
Set-CellText $ws 'C18' 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText $ws 'D18' '30.467.96'
Set-CellText $ws 'E18' '  -0.12%  '

Set-CellText $ws 'D19' '1.002'
Set-CellText $ws 'E19' '  +0.16%  '

Set-CellText $ws 'B20' 'Avalanche'
Set-CellText $ws 'C20' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-CellText $ws 'D20' '12.98'
Set-CellText $ws 'E20' '  +2.45%  '

Set-CellText $ws 'D21' '2.163.64'
Set-CellText $ws 'E21' '  +3.53%  '

Set-CellText $ws 'B22' 'ShibaInu'
Set-CellText $ws 'C22' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-CellText $ws 'D22' '0.000007494'
Set-CellText $ws 'E22' '  +2.31%  '

Set-CellText $ws 'B23' 'BinanceUSD'
Set-CellText $ws 'C23' 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
Set-CellText $ws 'D23' '1.001'
Set-CellText $ws 'E23' '  +0.02%  '

Set-CellText $ws 'B24' 'Uniswap'
Set-CellText $ws 'C24' 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-CellText $ws 'D24' '5.197'
Set-CellText $ws 'E24' '  +5.53%  '

Set-CellText $ws 'B25' 'Chainlink'
Set-CellText $ws 'C25' 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-CellText $ws 'D25' '6.338'
Set-CellText $ws 'E25' '  +7.07%  '

Set-CellText $ws 'B26' 'Cosmos'
Set-CellText $ws 'C26' 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-CellText $ws 'D26' '9.284'
Set-CellText $ws 'E26' '  +1.45%  '

Set-CellText $ws 'B27' 'Monero'
Set-CellText $ws 'C27' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-CellText $ws 'D27' '165.18'
Set-CellText $ws 'E27' '  +1.59%  '

Set-CellText $ws 'B28' 'EthereumClassic'
Set-CellText $ws 'C28' 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-CellText $ws 'D28' '19.85'
Set-CellText $ws 'E28' '  +10.56%  '

Set-CellText $ws 'B29' 'LidoDAOToken'
Set-CellText $ws 'C29' 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-CellText $ws 'D29' '2.037'
Set-CellText $ws 'E29' '  +8.66%  '

Set-CellText $ws 'B30' 'Stellar'
Set-CellText $ws 'C30' 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-CellText $ws 'D30' '0.1127'
Set-CellText $ws 'E30' '  +10.52%  '

Set-CellText $ws 'B31' 'Toncoin'
Set-CellText $ws 'C31' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-CellText $ws 'D31' '1.349'
Set-CellText $ws 'E31' '  -0.79%  '

Set-CellText $ws 'B32' 'InternetComputer(DFINITY)'
Set-CellText $ws 'C32' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 'D32' '4.128'
Set-CellText $ws 'E32' '  +3.15%  '

Set-CellText $ws 'B33' 'Filecoin'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D33' '3.939'
Set-CellText $ws 'E33' '  +3.02%  '

Set-CellText $ws 'B34' 'Hedera'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-CellText $ws 'D34' '0.05034'
Set-CellText $ws 'E34' '  +4.11%  '

Set-CellText $ws 'B35' 'ImmutableX'
Set-CellText $ws 'C35' 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-CellText $ws 'D35' '0.7388'
Set-CellText $ws 'E35' '  +5.44%  '

Set-CellText $ws 'B36' 'ARBITRUM'
Set-CellText $ws 'C36' 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-CellText $ws 'D36' '1.149'
Set-CellText $ws 'E36' '  +1.51%  '

Set-CellText $ws 'B37' 'Frax'
Set-CellText $ws 'C37' 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
Set-CellText $ws 'D37' '1.001'
Set-CellText $ws 'E37' '  +0.17%  '

Set-CellText $ws 'B38' 'HuobiToken'
Set-CellText $ws 'C38' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws 'D38' '2.718'
Set-CellText $ws 'E38' '  +0.53%  '

Set-CellText $ws 'B39' 'VeChain'
Set-CellText $ws 'C39' 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws 'D39' '0.01969'
Set-CellText $ws 'E39' '  +3.62%  '

Set-CellText $ws 'B40' 'MXToken'
Set-CellText $ws 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-CellText $ws 'D40' '2.704'
Set-CellText $ws 'E40' '  +1.09%  '

Set-CellText $ws 'B41' 'RenderToken'
Set-CellText $ws 'C41' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-CellText $ws 'D41' '2.024'
Set-CellText $ws 'E41' '  +1.83%  '

Set-CellText $ws 'B42' 'TrustWalletToken'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-CellText $ws 'D42' '0.8787'
Set-CellText $ws 'E42' '  +0.58%  '

Set-CellText $ws 'B43' 'Quant'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-CellText $ws 'D43' '107.06'
Set-CellText $ws 'E43' '  +0.33%  '

Set-CellText $ws 'B44' 'FraxShare'
Set-CellText $ws 'C44' 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-CellText $ws 'D44' '5.910'
Set-CellText $ws 'E44' '  +7.06%  '

Set-CellText $ws 'B45' 'PaxDollar'
Set-CellText $ws 'C45' 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
Set-CellText $ws 'D45' '1.001'
Set-CellText $ws 'E45' '  +0.14%  '

Set-CellText $ws 'B46' 'Aave'
Set-CellText $ws 'C46' 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-CellText $ws 'D46' '68.86'
Set-CellText $ws 'E46' '  +10.42%  '

Set-CellText $ws 'B47' 'TheSandbox'
Set-CellText $ws 'C47' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-CellText $ws 'D47' '0.4155'
Set-CellText $ws 'E47' '  +2.50%  '

Set-CellText $ws 'B48' 'Aptos'
Set-CellText $ws 'C48' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-CellText $ws 'D48' '7.259'
Set-CellText $ws 'E48' '  +1.22%  '

Set-CellText $ws 'B49' 'EnergySwap'
Set-CellText $ws 'C49' 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-CellText $ws 'D49' '9.310'
Set-CellText $ws 'E49' '  +8.79%  '

Set-CellText $ws 'B50' 'Algorand'
Set-CellText $ws 'C50' 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-CellText $ws 'D50' '0.1212'
Set-CellText $ws 'E50' '  +0.26%  '

Set-CellText $ws 'B51' 'Elrond'
Set-CellText $ws 'C51' 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-CellText $ws 'D51' '34.66'
Set-CellText $ws 'E51' '  +3.61%  '
